$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data rows (68 = site A, 69 = site B)
# down onto the two new rows so the new cells reuse the same style indices
# (no new numFmts / cellXfs get created) instead of Excel synthesizing new ones.
$ws.Range("A68:F68").Copy()
$ws.Range("A70:F70").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A69:F69").Copy()
$ws.Range("A71:F71").PasteSpecial(-4122)  # xlPasteFormats

# Row 70: 四方坪站
$ws.Range("A70").Value = 45965
$ws.Range("B70").Value = "四方坪站"
$ws.Range("C70").Value = 9285.74
$ws.Range("D70").Value = 8020.26
$ws.Range("E70").Value = 3136.57
$ws.Range("F70").Value = 393

# Row 71: 高岭站
$ws.Range("A71").Value = 45965
$ws.Range("B71").Value = "高岭站"
$ws.Range("C71").Value = 3837.03
$ws.Range("D71").Value = 3288.66
$ws.Range("E71").Value = 1047.0899999999999
$ws.Range("F71").Value = 152

# Update the active cell / selection to match the new end of data
$ws.Range("H67").Select()
